$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Step 1: Populate the new rows 22-30 for columns A, D (serial_num), then E (ip_address), then B (name), and F/G/H/I/J --
# matching the order in which the author appears to have entered the data
# (new shared strings appear in this column-major order in the sharedStrings table).

$newIds = @(10021, 10022, 10023, 10024, 10025, 10026, 10027, 10028, 10029)
$newSerials = @("FB5962911653", "FB5962911654", "FB5962911655", "FB5962911656", "FB5962911657", "FB5962911658", "FB5962911659", "FB5962911661", "FB5962911662")
$newIps = @("192.168.0.874", "192.168.0.721", "192.168.0.841", "192.168.0.186", "192.168.0.627", "192.168.0.879", "192.168.0.628", "192.168.0.306", "192.168.0.355")
$newNames = @("Machine 21", "Machine 22", "Machine 23", "Machine 24", "Machine 25", "Machine 26", "Machine 27", "Machine 28", "Machine 29")

for ($i = 0; $i -lt 9; $i++) {
    $r = 22 + $i
    $ws.Cells.Item($r, 1).Value = $newIds[$i]
}
for ($i = 0; $i -lt 9; $i++) {
    $r = 22 + $i
    $ws.Cells.Item($r, 4).Value = $newSerials[$i]
}
for ($i = 0; $i -lt 9; $i++) {
    $r = 22 + $i
    $ws.Cells.Item($r, 5).Value = $newIps[$i]
}
for ($i = 0; $i -lt 9; $i++) {
    $r = 22 + $i
    $ws.Cells.Item($r, 2).Value = $newNames[$i]
}
for ($i = 0; $i -lt 9; $i++) {
    $r = 22 + $i
    $ws.Cells.Item($r, 6).Value = 1001
    $ws.Cells.Item($r, 7).Value = "eng"
    $ws.Cells.Item($r, 8).Value = $true
    $ws.Cells.Item($r, 9).Value = "superadmin"
    $ws.Cells.Item($r, 10).Value = "now()"
}

# Step 2: Replace mac_address (column C) for ALL rows 2-30 with the new
# dash-separated, upper-case format -- this is done last so the new MAC
# strings land at the tail of the shared-strings table, after the row-22..30
# serial/ip/name strings, matching the authored diff.
$newMacs = @(
    "8C-16-45-5A-5D-0D",
    "8C-16-45-88-E1-0D",
    "00-FF-D3-E3-9A-27",
    "8C-16-45-5A-62-41",
    "E8-6A-64-1D-75-E4",
    "8C-16-45-FA-94-B7",
    "8C-16-45-1A-0F-62",
    "E8-6A-64-1C-52-6E",
    "48-51-B7-10-35-A6",
    "8C-16-45-38-F3-F3",
    "D4-3D-7E-58-CC-45",
    "8C-16-45-5A-5D-96",
    "8C-16-45-5A-5D-8E",
    "8C-16-45-33-A5-5F",
    "3C-95-09-F9-EA-DF",
    "8C-16-45-88-E7-0B",
    "B4-69-21-5A-DB-C4",
    "E8-6A-64-1D-48-B7",
    "8C-16-45-59-69-09 ",
    "98-E7-F4-30-16-5A ",
    "38-BA-F8-53-C7-8F",
    "E8-6A-64-1C-58-C2",
    "E4-A4-71-CE-BA-93",
    "54-E1-AD-EA-30-C9",
    "8C-16-45-65-DD-40",
    "58-20-B1-D6-C3-BE",
    "8C-16-45-38-F0-25",
    "6C-88-14-AC-EF-55",
    "3C-6A-A7-C0-DF-27"
)
for ($i = 0; $i -lt $newMacs.Length; $i++) {
    $r = 2 + $i
    $ws.Cells.Item($r, 3).Value = $newMacs[$i]
}

# Step 3: column C width grew slightly once the data changed (Excel's
# best-fit auto-width recalculated for the new content).
$ws.Columns("C").ColumnWidth = 16.14

# Step 4: the author ended by selecting from row 31 to the bottom of the sheet
# (e.g. selecting the empty rows below the table).
$ws.Range("A31:XFD1048576").Select()
